# feat: add 2022-Q3 data
# Inserts a new "2022-Q3" worksheet (fund holdings detail) ahead of "2022-Q2",
# and rolls the new quarter's totals into the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "2022-Q3" right before "2022-Q2" ---
$q2Ref = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Ref)
$newSheet.Name = "2022-Q3"
# Re-resolve a fresh handle to "2022-Q2" -- the handle used for positioning the
# Add() call above gets rebound to the freshly-inserted sheet, not the original.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- 2. Populate the new "2022-Q3" sheet with the fund holdings table ---
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# columns B, D, E, F, G hold numeric-looking data that is stored as *text* in
# the workbook (fund codes need their leading zeros kept, and the metrics
# columns match the text formatting used on every other quarter sheet)
$newSheet.Range("B2:B16").NumberFormat = "@"
$newSheet.Range("D2:G16").NumberFormat = "@"

$data = @(
    @("161017", "富国中证500指数增强（LOF）", "66.37", "90.18", "0.78", "0.5177", 9),
    @("004932", "招商丰拓灵活配置混合A", "5.17", "88.09", "5.79", "0.2993", 5),
    @("010673", "兴全中证800六个月持有期指数增强A", "12.02", "96.73", "2.41", "0.2897", 10),
    @("002657", "招商安裕灵活配置混合A", "16.92", "31.74", "1.65", "0.2792", 10),
    @("004933", "招商丰拓灵活配置混合C", "4.12", "88.09", "5.79", "0.2385", 5),
    @("002581", "招商丰凯灵活配置混合A", "4.22", "48.56", "2.44", "0.1030", 3),
    @("004143", "招商盛合灵活配置混合C", "3.37", "50.91", "2.83", "0.0954", 8),
    @("002658", "招商安裕灵活配置混合A", "4.40", "31.74", "1.65", "0.0726", 10),
    @("002510", "申万菱信中证500指数增强A", "3.74", "93.28", "1.86", "0.0696", 5),
    @("010674", "兴全中证800六个月持有期指数增强C", "1.29", "96.73", "2.41", "0.0311", 10),
    @("002582", "招商丰凯灵活配置混合C", "1.15", "48.56", "2.44", "0.0281", 3),
    @("015206", "招商安裕灵活配置混合D", "1.21", "31.74", "1.65", "0.0200", 10),
    @("007795", "申万菱信中证500指数增强C", "0.87", "93.28", "1.86", "0.0162", 5),
    @("004142", "招商盛合灵活配置混合A", "0.57", "50.91", "2.83", "0.0161", 8),
    @("013332", "富国中证500指数增强(LOF)C", "1.68", "90.18", "0.78", "0.0131", 9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $newSheet.Cells.Item($row, 1).Value = $i
    $newSheet.Cells.Item($row, 2).Value = $rec[0]
    $newSheet.Cells.Item($row, 3).Value = $rec[1]
    $newSheet.Cells.Item($row, 4).Value = $rec[2]
    $newSheet.Cells.Item($row, 5).Value = $rec[3]
    $newSheet.Cells.Item($row, 6).Value = $rec[4]
    $newSheet.Cells.Item($row, 7).Value = $rec[5]
    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}

# header row + index column styling (bold, centered, bordered) to match the
# other quarter sheets -- copy formats only from the "2022-Q2" sheet so the
# text NumberFormat set above is left untouched on the data columns
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2:A16").PasteSpecial(-4122)

# --- 3. Update the "总计" (summary) sheet: shift rows down and prepend 2022-Q3 totals ---
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @("2022-Q3", 15, 2.09),
    @("2022-Q2", 29, 3.06),
    @("2022-Q1", 28, 2.98),
    @("2021-Q4", 24, 2.88),
    @("2021-Q2", 4, 0.44),
    @("2021-Q1", 1, 0.06),
    @("2020-Q4", 2, 0.05)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $row = $i + 2
    $rec = $summaryRows[$i]
    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $rec[0]
    $summary.Cells.Item($row, 3).Value = $rec[1]
    $summary.Cells.Item($row, 4).Value = $rec[2]
}

# carry the bold/centered index style onto the newly-added A8 cell
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)
